$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 212 (pushes the old rows 212..331 down to 214..333).
$ws.Rows("212:213").Insert()

# Row 212 is the first freshly-inserted row; fill it in with its full record.
$ws.Range("A212").Value2 = 4
$ws.Range("B212").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C212").Value2 = "Los Lagos"
$ws.Range("D212").Value2 = 44830
$ws.Range("E212").Value2 = 10
$ws.Range("F212").Value2 = 100112003
$ws.Range("G212").Value2 = "Ajo"
$ws.Range("H212").Value2 = "Chino"
$ws.Range("I212").Value2 = "Primera"
$ws.Range("J212").Value2 = 50
$ws.Range("K212").Value2 = 23000
$ws.Range("L212").Value2 = 23000
$ws.Range("M212").Value2 = 23000
$ws.Range("N212").Value2 = "$/caja 10 kilos"
$ws.Range("O212").Value2 = "China"
$ws.Range("P212").Value2 = 2300
$ws.Range("Q212").Value2 = 10
$ws.Range("R212").Value2 = "Hortaliza"

# Row 213 is the second freshly-inserted row; fill it in with its full record.
$ws.Range("A213").Value2 = 4
$ws.Range("B213").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C213").Value2 = "Los Lagos"
$ws.Range("D213").Value2 = 44830
$ws.Range("E213").Value2 = 10
$ws.Range("F213").Value2 = 100112003
$ws.Range("G213").Value2 = "Ajo"
$ws.Range("H213").Value2 = "Chino"
$ws.Range("I213").Value2 = "Primera"
$ws.Range("J213").Value2 = 30
$ws.Range("K213").Value2 = 25000
$ws.Range("L213").Value2 = 25000
$ws.Range("M213").Value2 = 25000
$ws.Range("N213").Value2 = "$/malla 10 kilos"
$ws.Range("O213").Value2 = "China"
$ws.Range("P213").Value2 = 2500
$ws.Range("Q213").Value2 = 10
$ws.Range("R213").Value2 = "Hortaliza"

# Make sure the date columns keep the date style/format used elsewhere in column D.
$ws.Range("D212").NumberFormat = $ws.Range("D215").NumberFormat
$ws.Range("D213").NumberFormat = $ws.Range("D215").NumberFormat
